# ModeloCreatorAuthor.xlsx - add three new activity rows (Arrasta e Solta,
# Associar, Jogo da memoria) to the activity list, extend the formatted
# range down to row 33, drop the stray leftover formatted cell that used
# to live at D11, and refresh the column C width / selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old leftover formatted row (used to hold just D11) --------
$ws.Rows.Item(11).Delete()

# --- Propagate row 5's formatting down through row 33 so every new row --
# --- lines up with the existing "Tipo / Enunciado / Resposta / Opcoes" --
# --- table look (white cells, Arial 16, thin borders). -------------------
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Rows 2-5 all carry a 20.25pt row height (auto-expanded for the 16pt
# Arial font); give the newly formatted rows the same height.
$ws.Range("A6:G33").RowHeight = 20.25

# --- Row 6: "Arrasta e Solta" (drag and drop) activity -------------------
$ws.Range("A6").Value = "Arrasta e Solta"
$ws.Range("B6").Value = "Como é feito o chocolate?"
$ws.Range("C6").Value = "Cacau"
$ws.Range("D6").Value = "Morango"
$ws.Range("E6").Value = "Cacau"
$ws.Range("F6").Value = "Abacaxi"
$ws.Range("G6").Value = "Pera"

# --- Row 7: "Associar" (match) activity -----------------------------------
$ws.Range("A7").Value = "Associar"
$ws.Range("B7").Value = "Associe"
$ws.Range("C7").Value = "opção 1 - opção 3 - opção 2 - opção 4"
$ws.Range("D7").Value = "Carro"
$ws.Range("E7").Value = "Maçã"
$ws.Range("F7").Value = "Motor"
$ws.Range("G7").Value = "Arvore"

# --- Row 8: "Jogo da memória" (memory game) activity ----------------------
$ws.Range("A8").Value = "Jogo da memória"
$ws.Range("B8").Value = "Ache as cartas certas"
$ws.Range("C8").Value = "opção 1 - opção 3 - opção 2 - opção 4"
$ws.Range("D8").Value = "2 + 2"
$ws.Range("E8").Value = "4 + 4"
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 8

# Rows 9-33 stay blank (pre-formatted placeholder rows for future
# activities), already picked up the formatting from the paste above.

# --- Column C now holds the longest text in the sheet ("opção 1 - opção -
# --- 3 - opção 2 - opção 4"), so it needs to be widened / re-bestFit. -----
$ws.Columns.Item(3).AutoFit()

# --- Selection moved on to F9 (first blank placeholder row) when the -----
# --- author finished typing the new data. ---------------------------------
$ws.Range("F9").Select()
